# Normalize header / concept casing and add the missing "spoon" row label.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "apple_feature"
$ws.Range("C2").Value = "apple_val"
$ws.Range("D2").Value = "bowl_feature"
$ws.Range("E2").Value = "bowl_val"

$ws.Range("A3").Value = "apple"
$ws.Range("A4").Value = "bowl"
$ws.Range("A11").Value = "spoon"

# Restore the view to show column A with A12 selected (was scrolled to C1 / R12:S12).
$ws.Range("A12").Select()
